# DC-Colos.xlsx update
#
# The underlying data got re-sorted/re-generated upstream: within two
# contiguous blocks of rows the last row moved up to become the new first
# row of the block (everything else shifts down by one row). We reproduce
# that by reading each block into a 2-D array, rotating it, and writing it
# back - this avoids hand-transcribing ~60 rows x 7 columns of data.
#
# Block 1: rows 203-211 (A:G) - South America colo list
# Block 2: rows 235-284 (A:G) - North America colo list

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Rotate-Block($ws, $firstRow, $lastRow) {
    $addr = "A" + $firstRow + ":G" + $lastRow
    $rng = $ws.Range($addr)
    $vals = $rng.Value2

    $rows = $vals.GetLength(0)
    $cols = $vals.GetLength(1)

    $newvals = New-Object 'object[,]' $rows, $cols

    # new first row = old last row
    for ($c = 1; $c -le $cols; $c++) {
        $newvals[0, $c - 1] = $vals[$rows, $c]
    }
    # everything else shifts down by one
    for ($r = 2; $r -le $rows; $r++) {
        for ($c = 1; $c -le $cols; $c++) {
            $newvals[$r - 1, $c - 1] = $vals[$r - 1, $c]
        }
    }

    $rng.Value2 = $newvals
}

Rotate-Block $ws 203 211
Rotate-Block $ws 235 284

# The newly-promoted first row of block 1 (SJK / Sao Jose dos Campos) also
# gets its "name" column normalised to include the country, matching the
# "City, Country" convention used by every other row in the sheet.
$ws.Range("B203").Value = "São José dos Campos, Brazil"
